$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 16.45880000000001
$ws.Range("C3").Value = -12.6352
$ws.Range("E3").Value = 16.7851
$ws.Range("C4").Value = -11.96479999999999
$ws.Range("E9").Value = 17.24820000000001
$ws.Range("B11").Value = 5.853900000000001
$ws.Range("B12").Value = 5.256499999999999
$ws.Range("C14").Value = -13.24619999999999
$ws.Range("B15").Value = 4.690799999999999
$ws.Range("E15").Value = 16.2191
$ws.Range("E19").Value = 16.56080000000001
$ws.Range("E20").Value = 15.74779999999999
$ws.Range("E25").Value = 17.13860000000001
$ws.Range("C26").Value = -12.36880000000001
$ws.Range("B27").Value = 6.908600000000004
$ws.Range("E27").Value = 16.7853
$ws.Range("B28").Value = 6.816700000000005
$ws.Range("E28").Value = 16.7356
$ws.Range("E30").Value = 15.9265
$ws.Range("B31").Value = 4.046799999999997
$ws.Range("C31").Value = -13.4781
$ws.Range("B32").Value = 6.645099999999998
$ws.Range("E32").Value = 16.61819999999998
$ws.Range("C35").Value = -13.22960000000001
$ws.Range("B36").Value = 9.116100000000003
$ws.Range("C37").Value = -13.63869999999999
$ws.Range("B38").Value = 5.829199999999999
$ws.Range("C39").Value = -12.0332
$ws.Range("C40").Value = -13.23910000000001
$ws.Range("E44").Value = 16.8459
$ws.Range("C45").Value = -13.12239999999999
$ws.Range("B46").Value = 7.616299999999998
$ws.Range("E47").Value = 16.7605
$ws.Range("C52").Value = -10.65969999999999
$ws.Range("B54").Value = 4.6287
$ws.Range("B55").Value = 4.823399999999999
$ws.Range("B56").Value = 4.541399999999999
$ws.Range("C57").Value = -14.81089999999999
$ws.Range("E58").Value = 16.34820000000001
$ws.Range("E62").Value = 16.66370000000001
$ws.Range("B67").Value = 6.079399999999997
$ws.Range("B69").Value = 5.794799999999996
$ws.Range("B72").Value = 5.623800000000005
$ws.Range("B73").Value = 8.549500000000002
$ws.Range("E77").Value = 17.22140000000002
$ws.Range("E78").Value = 16.47080000000002
$ws.Range("C81").Value = -13.2639
$ws.Range("B83").Value = 5.153699999999994
$ws.Range("C83").Value = -14.2456
$ws.Range("E84").Value = 16.8347
$ws.Range("B86").Value = 4.793200000000002
$ws.Range("E89").Value = 17.22810000000001
$ws.Range("B91").Value = 5.3699
$ws.Range("E91").Value = 17.85470000000002
$ws.Range("E92").Value = 18.04070000000003
$ws.Range("B93").Value = 6.987499999999999
$ws.Range("E96").Value = 16.78809999999998
$ws.Range("B99").Value = 4.826199999999997
$ws.Range("C100").Value = -13.78059999999999
$ws.Range("C102").Value = -11.9945
$ws.Range("E102").Value = 16.90669999999999
